$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update item descriptions (column E) to remove duplicate/erroneous line items
$ws.Cells.Item(4, 5).Value = "SSD1TB - 1 - SSD SATA 2.5  1TB - 850.0`nSERV006 - 1 - Bat He ZBook 15 G1 AR08 Negra 8C - 359`nENVIOS - 1 - Envío - 72.5"
$ws.Cells.Item(5, 5).Value = "SERV006 - 1 - EKD6600  Bateria Li-ion 11.1V para Dell Precision M4600, M4700, M6600 - 365.4`nENVIOS - 1 - Envío - 50.73"
$ws.Cells.Item(10, 5).Value = "SERV03 - 1 - Mantenimiento general - 250.0`nSERV006 - 1 - XEKHHT03XL Bateria Interna (P) 11.4V para HP Pavilion 14-CE 14-CF 14-cm 14-CK 14-DF 14-MA 15-DA 15-CS 15-DB 14Q-CS 15-CW 17-by 17-CA HP 240 245 250 255 G7 340 348 G5 HT03XL - 324.8`nENVIOS - 1 - Envío - 50.73"
$ws.Cells.Item(28, 5).Value = "SERV006 - 1 - Bat To Sat L50-B PA5184U-1BRS Negra 4C - 229`nENVIOS - 1 - Envío - 72.5"
$ws.Cells.Item(41, 5).Value = "REF COMODÍN - 1 - TEC387 Teclado color Negro (SP) para IBM Lenovo Thinkpad T530 X230 - 481.4`nENVIOS - 1 - Envío - 7.25"
$ws.Cells.Item(60, 5).Value = "REF COMODÍN - 1 - Cargador Asus Vivobook X441 X556 19v 3.42a 65w 4.0*1.35mm - 275"
$ws.Cells.Item(63, 5).Value = "REF COMODÍN - 1 - EKD6MT4T Bateria Interna (P) 7.6V para Dell Latitude E5470/ E5570 Precision 3510 HK6DV 079VRK TXF9M 0TXF9M 6MT4T - 446.31`nENVIOS - 1 - Envío - 68.09"
$ws.Cells.Item(82, 5).Value = "REF COMODÍN - 1 - EKDGK5KY Bateria 11.4V 43Wh para Dell Inspiron 11 3000 series / 13 7000 series / 15 7000 GK5KY - 440.8`nENVIOS - 1 - Envío - 68.09"
$ws.Cells.Item(92, 5).Value = "REF COMODÍN - 1 - TEC491 Teclado color Negro (SP) para Lenovo E431 T440 L440 E440 T450 T460 - 336.4`nREF COMODÍN - 1 - TEC526 Teclado color Negro (SP) para Lenovo E531 E540 L540 L570 T540 T550 W540 W550s Series - 452.4`nREF COMODÍN - 1 - LCD140-015D Pantalla LCD 14.0 LED WXGA (1366X768) Slim Conector Inferior Derecho 30P Glossy 315mm (bezel delgado) - 696`nENVIOS - 1 - Envío - 20.66"
$ws.Cells.Item(95, 5).Value = "SSD256 - 1 - SSD SATA 2.5  256GB - 650.86`nREF COMODÍN - 1 - Bateria para ACER BT12122 - 245"

# Update corresponding cost totals (column F) to reflect removed line items
$ws.Cells.Item(4, 6).Value = 1281.5
$ws.Cells.Item(5, 6).Value = 416.13
$ws.Cells.Item(10, 6).Value = 625.53
$ws.Cells.Item(28, 6).Value = 301.5
$ws.Cells.Item(41, 6).Value = 488.65
$ws.Cells.Item(60, 6).Value = 275
$ws.Cells.Item(63, 6).Value = 514.4
$ws.Cells.Item(82, 6).Value = 508.89
$ws.Cells.Item(92, 6).Value = 1505.46
$ws.Cells.Item(95, 6).Value = 895.86

# Add totals row with SUM formulas
$ws.Cells.Item(135, 6).Formula = "=SUM(F2:F134)"
$ws.Cells.Item(135, 7).Formula = "=SUM(G2:G134)"
$ws.Cells.Item(135, 8).Formula = "=G135-F135"
